# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45221 (2023-10-22) to 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 146

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
